$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Amir Temur question
$ws.Range("A2").Value = "Amir Temur qachon tug‘ilgan?"
$ws.Range("B2").Value = "'1336"
$ws.Range("C2").Value = "'1396"
$ws.Range("D2").Value = "'1405"
$ws.Range("E2").Value = "A"
$ws.Range("B2:D2").Style = "Normal"

# Row 3: Buxoro question
$ws.Range("A3").Value = "Buxoro qayerda joylashgan?"
$ws.Range("B3").Value = "Toshkent"
$ws.Range("C3").Value = "Buxoro"
$ws.Range("D3").Value = "Samarqand"
$ws.Range("E3").Value = "B"

# Remove old row 4 entirely - table shrinks to 2 data rows
$ws.Range("A4:E4").Delete()
